# Auto-generated Excel COM-interop script
# Updates loading_percent values in Sheet1 (columns B,C,E,F,G,I,K,L,M for rows 2-25)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 15.56215833890737
$ws.Cells.Item(2, 3).Value = 5.759714878160847
$ws.Cells.Item(2, 5).Value = 9.523874639057887
$ws.Cells.Item(2, 6).Value = 16.86991607391233
$ws.Cells.Item(2, 7).Value = 3.781123445088729
$ws.Cells.Item(2, 9).Value = 42.16219757563169
$ws.Cells.Item(2, 11).Value = 14.40248692811434
$ws.Cells.Item(2, 12).Value = 10.70372182804004
$ws.Cells.Item(2, 13).Value = 16.24281408149997
$ws.Cells.Item(3, 2).Value = 15.52568906538251
$ws.Cells.Item(3, 3).Value = 5.563896747641673
$ws.Cells.Item(3, 5).Value = 9.540716599497955
$ws.Cells.Item(3, 6).Value = 15.89584955866808
$ws.Cells.Item(3, 7).Value = 3.784237520269339
$ws.Cells.Item(3, 9).Value = 41.80304197212396
$ws.Cells.Item(3, 11).Value = 14.36456919071079
$ws.Cells.Item(3, 12).Value = 10.70974789075377
$ws.Cells.Item(3, 13).Value = 16.26228532637093
$ws.Cells.Item(4, 2).Value = 15.50825481619839
$ws.Cells.Item(4, 3).Value = 5.44140202283918
$ws.Cells.Item(4, 5).Value = 9.551935103544293
$ws.Cells.Item(4, 6).Value = 15.26997757108491
$ws.Cells.Item(4, 7).Value = 3.786247852858394
$ws.Cells.Item(4, 9).Value = 41.58268005752307
$ws.Cells.Item(4, 11).Value = 14.34519809636557
$ws.Cells.Item(4, 12).Value = 10.71486821830227
$ws.Cells.Item(4, 13).Value = 16.27732367590934
$ws.Cells.Item(5, 2).Value = 15.50240363487425
$ws.Cells.Item(5, 3).Value = 5.390995233037257
$ws.Cells.Item(5, 5).Value = 9.556727826232748
$ws.Cells.Item(5, 6).Value = 15.00819731993403
$ws.Cells.Item(5, 7).Value = 3.787091882266368
$ws.Cells.Item(5, 9).Value = 41.49296711216767
$ws.Cells.Item(5, 11).Value = 14.33829375560811
$ws.Cells.Item(5, 12).Value = 10.71731225387776
$ws.Cells.Item(5, 13).Value = 16.28422720587675
$ws.Cells.Item(6, 2).Value = 15.50150792031142
$ws.Cells.Item(6, 3).Value = 5.382598126988011
$ws.Cells.Item(6, 5).Value = 9.557537020571779
$ws.Cells.Item(6, 6).Value = 14.96433081551593
$ws.Cells.Item(6, 7).Value = 3.787233533357639
$ws.Cells.Item(6, 9).Value = 41.47807687938717
$ws.Cells.Item(6, 11).Value = 14.33720722020448
$ws.Cells.Item(6, 12).Value = 10.71773968159463
$ws.Cells.Item(6, 13).Value = 16.28542035456953
$ws.Cells.Item(7, 2).Value = 15.50817082239522
$ws.Cells.Item(7, 3).Value = 5.440724090180858
$ws.Cells.Item(7, 5).Value = 9.551998844092173
$ws.Cells.Item(7, 6).Value = 15.26647399323137
$ws.Cells.Item(7, 7).Value = 3.786259135191605
$ws.Cells.Item(7, 9).Value = 41.58146974815751
$ws.Cells.Item(7, 11).Value = 14.34510096807761
$ws.Cells.Item(7, 12).Value = 10.71489973171747
$ws.Cells.Item(7, 13).Value = 16.27741364034286
$ws.Cells.Item(8, 2).Value = 15.54855914345538
$ws.Cells.Item(8, 3).Value = 5.692711875996853
$ws.Cells.Item(8, 5).Value = 9.529499894260537
$ws.Cells.Item(8, 6).Value = 16.53996406344768
$ws.Cells.Item(8, 7).Value = 3.782176834813879
$ws.Cells.Item(8, 9).Value = 42.03835503001746
$ws.Cells.Item(8, 11).Value = 14.38860530518333
$ws.Cells.Item(8, 12).Value = 10.70550503904622
$ws.Cells.Item(8, 13).Value = 16.24888799084592
$ws.Cells.Item(9, 2).Value = 15.66673902899761
$ws.Cells.Item(9, 3).Value = 6.165753325344673
$ws.Cells.Item(9, 5).Value = 9.492322293766744
$ws.Cells.Item(9, 6).Value = 19.00274580682531
$ws.Cells.Item(9, 7).Value = 3.774947145068249
$ws.Cells.Item(9, 9).Value = 42.93354423590964
$ws.Cells.Item(9, 11).Value = 14.5046237418302
$ws.Cells.Item(9, 12).Value = 10.69833599262844
$ws.Cells.Item(9, 13).Value = 16.21740256726471
$ws.Cells.Item(10, 2).Value = 15.7767476919666
$ws.Cells.Item(10, 3).Value = 6.496661936300721
$ws.Cells.Item(10, 5).Value = 9.469213959903653
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 3.770102624536171
$ws.Cells.Item(10, 9).Value = 43.58801845744056
$ws.Cells.Item(10, 11).Value = 14.60808048073713
$ws.Cells.Item(10, 12).Value = 10.69990613939928
$ws.Cells.Item(10, 13).Value = 16.20915729182157
$ws.Cells.Item(11, 2).Value = 15.83167419810031
$ws.Cells.Item(11, 3).Value = 6.642920873384446
$ws.Cells.Item(11, 5).Value = 9.459609080022878
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 3.767998933080531
$ws.Cells.Item(11, 9).Value = 43.8844624398632
$ws.Cells.Item(11, 11).Value = 14.65897294384466
$ws.Cells.Item(11, 12).Value = 10.70209854502106
$ws.Cells.Item(11, 13).Value = 16.20863034939318
$ws.Cells.Item(12, 2).Value = 15.85316066814612
$ws.Cells.Item(12, 3).Value = 6.697643339235047
$ws.Cells.Item(12, 5).Value = 9.456101970119171
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 3.767216621089586
$ws.Cells.Item(12, 9).Value = 43.99648345183591
$ws.Cells.Item(12, 11).Value = 14.67878288630216
$ws.Cells.Item(12, 12).Value = 10.70314057864852
$ws.Cells.Item(12, 13).Value = 16.20889325740747
$ws.Cells.Item(13, 2).Value = 15.84850287291997
$ws.Cells.Item(13, 3).Value = 6.685888054567865
$ws.Cells.Item(13, 5).Value = 9.45685151147884
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 3.767384470863697
$ws.Cells.Item(13, 9).Value = 43.97236903787939
$ws.Cells.Item(13, 11).Value = 14.6744927402455
$ws.Cells.Item(13, 12).Value = 10.70290675010211
$ws.Cells.Item(13, 13).Value = 16.20881608849018
$ws.Cells.Item(14, 2).Value = 15.83342822587104
$ws.Cells.Item(14, 3).Value = 6.647436459361476
$ws.Cells.Item(14, 5).Value = 9.459317943922622
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 3.767934285486298
$ws.Cells.Item(14, 9).Value = 43.89368350004111
$ws.Cells.Item(14, 11).Value = 14.66059201142058
$ws.Cells.Item(14, 12).Value = 10.7021800331547
$ws.Cells.Item(14, 13).Value = 16.20864271816329
$ws.Cells.Item(15, 2).Value = 15.8242835645674
$ws.Cells.Item(15, 3).Value = 6.623796121200784
$ws.Cells.Item(15, 5).Value = 9.460845629235331
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 3.768272923935575
$ws.Cells.Item(15, 9).Value = 43.84545407949118
$ws.Cells.Item(15, 11).Value = 14.65214709671498
$ws.Cells.Item(15, 12).Value = 10.70176245952684
$ws.Cells.Item(15, 13).Value = 16.20859671066302
$ws.Cells.Item(16, 2).Value = 15.77325510558395
$ws.Cells.Item(16, 3).Value = 6.487013339605263
$ws.Cells.Item(16, 5).Value = 9.469859882123524
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 3.770242112843071
$ws.Cells.Item(16, 9).Value = 43.5686149582038
$ws.Cells.Item(16, 11).Value = 14.6048305412385
$ws.Cells.Item(16, 12).Value = 10.69979254341763
$ws.Cells.Item(16, 13).Value = 16.20925650920621
$ws.Cells.Item(17, 2).Value = 15.74319105833807
$ws.Cells.Item(17, 3).Value = 6.401971021347406
$ws.Cells.Item(17, 5).Value = 9.475621904717665
$ws.Cells.Item(17, 6).Value = 20.20408069617459
$ws.Cells.Item(17, 7).Value = 3.771475725397421
$ws.Cells.Item(17, 9).Value = 43.39841842657862
$ws.Cells.Item(17, 11).Value = 14.57677551674591
$ws.Cells.Item(17, 12).Value = 10.69896215334456
$ws.Cells.Item(17, 13).Value = 16.21048624824705
$ws.Cells.Item(18, 2).Value = 15.72635961671603
$ws.Cells.Item(18, 3).Value = 6.352657233921242
$ws.Cells.Item(18, 5).Value = 9.479021483636728
$ws.Cells.Item(18, 6).Value = 19.95656407809808
$ws.Cells.Item(18, 7).Value = 3.772194694128181
$ws.Cells.Item(18, 9).Value = 43.30040853978291
$ws.Cells.Item(18, 11).Value = 14.56100027204971
$ws.Cells.Item(18, 12).Value = 10.69862371392907
$ws.Cells.Item(18, 13).Value = 16.21149709080827
$ws.Cells.Item(19, 2).Value = 15.72074033058358
$ws.Cells.Item(19, 3).Value = 6.335893345609459
$ws.Cells.Item(19, 5).Value = 9.480187205311342
$ws.Cells.Item(19, 6).Value = 19.87204792380562
$ws.Cells.Item(19, 7).Value = 3.772439746341604
$ws.Cells.Item(19, 9).Value = 43.26720550542918
$ws.Cells.Item(19, 11).Value = 14.55572146129582
$ws.Cells.Item(19, 12).Value = 10.69853304931038
$ws.Cells.Item(19, 13).Value = 16.21189150707161
$ws.Cells.Item(20, 2).Value = 15.74634386288641
$ws.Cells.Item(20, 3).Value = 6.411065655041121
$ws.Cells.Item(20, 5).Value = 9.474999690301036
$ws.Cells.Item(20, 6).Value = 20.24955283636157
$ws.Cells.Item(20, 7).Value = 3.771343430150859
$ws.Cells.Item(20, 9).Value = 43.41654865895499
$ws.Cells.Item(20, 11).Value = 14.57972472105154
$ws.Cells.Item(20, 12).Value = 10.69903615036965
$ws.Cells.Item(20, 13).Value = 16.21032393389878
$ws.Cells.Item(21, 2).Value = 15.83783749738995
$ws.Cells.Item(21, 3).Value = 6.658748964231831
$ws.Cells.Item(21, 5).Value = 9.458589966661128
$ws.Cells.Item(21, 6).Value = 21.46857628470577
$ws.Cells.Item(21, 7).Value = 3.767772403899179
$ws.Cells.Item(21, 9).Value = 43.91680216440046
$ws.Cells.Item(21, 11).Value = 14.66466049430293
$ws.Cells.Item(21, 12).Value = 10.70238774518732
$ws.Cells.Item(21, 13).Value = 16.2086811007065
$ws.Cells.Item(22, 2).Value = 15.90163030873968
$ws.Cells.Item(22, 3).Value = 6.816743611331209
$ws.Cells.Item(22, 5).Value = 9.448623103656908
$ws.Cells.Item(22, 6).Value = 22.22866616901552
$ws.Cells.Item(22, 7).Value = 3.765521904230718
$ws.Cells.Item(22, 9).Value = 44.24234790539611
$ws.Cells.Item(22, 11).Value = 14.72330074548025
$ws.Cells.Item(22, 12).Value = 10.7058124353548
$ws.Cells.Item(22, 13).Value = 16.21030223269505
$ws.Cells.Item(23, 2).Value = 15.86722254088271
$ws.Cells.Item(23, 3).Value = 6.73278845819014
$ws.Cells.Item(23, 5).Value = 9.45387339808965
$ws.Cells.Item(23, 6).Value = 21.82633154458858
$ws.Cells.Item(23, 7).Value = 3.766715437552091
$ws.Cells.Item(23, 9).Value = 44.06874309342845
$ws.Cells.Item(23, 11).Value = 14.69172135823665
$ws.Cells.Item(23, 12).Value = 10.7038719477798
$ws.Cells.Item(23, 13).Value = 16.20919086429214
$ws.Cells.Item(24, 2).Value = 15.74491706845275
$ws.Cells.Item(24, 3).Value = 6.406955281347814
$ws.Cells.Item(24, 5).Value = 9.475280722686975
$ws.Cells.Item(24, 6).Value = 20.22900810905287
$ws.Cells.Item(24, 7).Value = 3.771403210459774
$ws.Cells.Item(24, 9).Value = 43.40835248035025
$ws.Cells.Item(24, 11).Value = 14.57839028228407
$ws.Cells.Item(24, 12).Value = 10.69900226341361
$ws.Cells.Item(24, 13).Value = 16.21039636967629
$ws.Cells.Item(25, 2).Value = 15.63064959237515
$ws.Cells.Item(25, 3).Value = 6.040456961671252
$ws.Cells.Item(25, 5).Value = 9.501639283045364
$ws.Cells.Item(25, 6).Value = 18.34778573295695
$ws.Cells.Item(25, 7).Value = 3.776820520717619
$ws.Cells.Item(25, 9).Value = 42.69173601903319
$ws.Cells.Item(25, 11).Value = 14.46999790313958
$ws.Cells.Item(25, 12).Value = 10.70372182804004
$ws.Cells.Item(25, 13).Value = 16.22330356202032
